# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new row 153 and 154) above the
# existing data block, pushing every following row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 153 (each Insert() pushes row 153.. down by one)
$ws.Rows.Item(153).Insert()
$ws.Rows.Item(153).Insert()

# New row 153 - "Primera" quality entry for 2022-08-25 (serial 44798)
$ws.Cells.Item(153, 1).Value = 7
$ws.Cells.Item(153, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(153, 3).Value = "Ñuble"
$ws.Cells.Item(153, 4).Value = 44798
$ws.Cells.Item(153, 5).Value = 16
$ws.Cells.Item(153, 6).Value = "Fruta"
$ws.Cells.Item(153, 7).Value = 100101
$ws.Cells.Item(153, 8).Value = "Berries"
$ws.Cells.Item(153, 9).Value = 100101007
$ws.Cells.Item(153, 10).Value = "Kiwi"
$ws.Cells.Item(153, 11).Value = "Hayward"
$ws.Cells.Item(153, 12).Value = "Primera"
$ws.Cells.Item(153, 13).Value = 120
$ws.Cells.Item(153, 14).Value = 7000
$ws.Cells.Item(153, 15).Value = 7500
$ws.Cells.Item(153, 16).Value = 7250
$ws.Cells.Item(153, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(153, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(153, 19).Value = 403
$ws.Cells.Item(153, 20).Value = 18

# New row 154 - "Segunda" quality entry for 2022-08-25 (serial 44798)
$ws.Cells.Item(154, 1).Value = 7
$ws.Cells.Item(154, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(154, 3).Value = "Ñuble"
$ws.Cells.Item(154, 4).Value = 44798
$ws.Cells.Item(154, 5).Value = 16
$ws.Cells.Item(154, 6).Value = "Fruta"
$ws.Cells.Item(154, 7).Value = 100101
$ws.Cells.Item(154, 8).Value = "Berries"
$ws.Cells.Item(154, 9).Value = 100101007
$ws.Cells.Item(154, 10).Value = "Kiwi"
$ws.Cells.Item(154, 11).Value = "Hayward"
$ws.Cells.Item(154, 12).Value = "Segunda"
$ws.Cells.Item(154, 13).Value = 60
$ws.Cells.Item(154, 14).Value = 6000
$ws.Cells.Item(154, 15).Value = 6000
$ws.Cells.Item(154, 16).Value = 6000
$ws.Cells.Item(154, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(154, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(154, 19).Value = 333
$ws.Cells.Item(154, 20).Value = 18
